# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new blank column immediately
# before the existing "Late" column (column N). This shifts the old N / O / P
# columns (Late / Date("heading") / Outstanding) one place to the right
# (-> O / P / Q) while leaving a new, empty column N in their place - exactly
# what Excel's "Insert Column" does when column N is selected and a column is
# inserted to its left.
#
# Then make "Repayment schedule" the active sheet/tab (it was "Transactions"
# before) and leave the selection on cell R8 (just past the last populated
# column/row of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N - shifts N->O, O->P, P->Q and carries over
# the cell styles from the old column into the new one (Excel's default
# "insert column" behaviour), matching the row/cell styles in the target.
$ws.Columns("N").Insert()

# The inserted column picks up the width Excel would copy from the column
# immediately to its left (column M) when doing an "Insert Column".
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab (was "Transactions"), and
# move the selection to cell R8.
$ws.Activate()
$ws.Range("R8").Select() | Out-Null
